$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    @("74-26=", "39+23="),
    @("47+25=", "16+77="),
    @("54-46=", "66-38="),
    @("43-25=", "6+75="),
    @("19+13=", "97-39="),
    @("84-68=", "70-14="),
    @("17+35=", "59+3="),
    @("39+18=", "18+78="),
    @("84-27=", "43-19="),
    @("49+45=", "54-7="),
    @("73-15=", "9+29="),
    @("18+14=", "35+57="),
    @("91-29=", "6+7="),
    @("28-9=", "18+46="),
    @("71-65=", "64+28="),
    @("83-39=", "28+27="),
    @("44-29=", "87+8="),
    @("82-79=", "7+47="),
    @("16+15=", "17+68="),
    @("47+15=", "11-5="),
    @("26+19=", "72-19="),
    @("25+38=", "28+44="),
    @("57-18=", "71-39="),
    @("75-67=", "43+28="),
    @("26+15=", "91-83="),
    @("12+49=", "74+19="),
    @("61-55=", "39+27="),
    @("6+37=", "71-63="),
    @("64-19=", "56+16="),
    @("31-8=", "27+57="),
    @("17+44=", "69+16="),
    @("66+18=", "62-28="),
    @("43-9=", "43-27="),
    @("27-19=", "23-19="),
    @("33+59=", "60-6="),
    @("42+9=", "16+38="),
    @("16+57=", "81-5="),
    @("9+5=", "90-77="),
    @("19+65=", "59+22="),
    @("19+34=", "44+28="),
    @("69+9=", "78+7="),
    @("36+15=", "53-7="),
    @("22-9=", "60-13="),
    @("91-6=", "40-19="),
    @("28+24=", "58+14="),
    @("35-7=", "91-42="),
    @("34-8=", "36+45="),
    @("66+6=", "7+74="),
    @("29+63=", "6+45="),
    @("35-7=", "62-58="),
    @("40-3=", "24+48="),
    @("95-9=", "7+29="),
    @("38+46=", "51-27="),
    @("19+42=", "16+25="),
    @("67+8=", "54-8="),
    @("35-26=", "38+53="),
    @("91-24=", "51-19="),
    @("19+49=", "80-15="),
    @("65-39=", "58+24="),
    @("9+55=", "37+8="),
    @("54+39=", "46-7="),
    @("83-49=", "82-9="),
    @("41-38=", "85+7="),
    @("93-47=", "50-48="),
    @("71-43=", "95-76="),
    @("82-14=", "95-89="),
    @("92-18=", "64-25="),
    @("27+65=", "29+3="),
    @("56-17=", "60-25="),
    @("40-7=", "39+4="),
    @("84-58=", "65+29="),
    @("75-47=", "6+89="),
    @("72-66=", "7+84="),
    @("27+68=", "87-69="),
    @("80-23=", "92-89="),
    @("33-8=", "6+68="),
    @("58+3=", "70-21="),
    @("18+3=", "67+26="),
    @("91-37=", "93-74="),
    @("28+39=", "40-14="),
    @("39+46=", "38+5="),
    @("42-5=", "32+49="),
    @("5+59=", "94-7="),
    @("68+25=", "71-18="),
    @("60-5=", "9+84="),
    @("58+25=", "17+18="),
    @("17+69=", "9+82="),
    @("50-38=", "12-3="),
    @("48+19=", "19+3="),
    @("72-37=", "17+9="),
    @("52-38=", "50-39="),
    @("33+9=", "77+17="),
    @("82-64=", "44+49="),
    @("64+17=", "30-3="),
    @("78-29=", "12-5="),
    @("15-6=", "9+22="),
    @("33+8=", "6+28="),
    @("60-57=", "94-79="),
    @("60-1=", "33-18="),
    @("92-55=", "96-77=")
)

$cols = 5
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = [int][math]::Floor($i / $cols) + 1
    $col = ($i % $cols) + 1
    $old = $values[$i][0]
    $new = $values[$i][1]
    $cell = $t.Cell($row, $col)
    $para = $cell.Range.Paragraphs.Item(1)
    $pr = $para.Range
    $cur = $pr.Text.TrimEnd([char]13, [char]7)
    if ($cur -ne $old) {
        Write-Host "MISMATCH at" $row $col "expected" $old "got" $cur
    }
    $pr.Text = $new
}
Write-Host "Done"
